# Applies the scheduled-runner profit/price recalculation updates to the
# leve-profit tables on each crafting-class worksheet (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 322558.78
$ws.Range("I64").Value = 540419.4399999999
$ws.Range("J64").Value = 4147.077
$ws.Range("K64").Value = 540419.4399999999
$ws.Range("L64").Value = 4147.077
$ws.Range("M64").Value = -540171.4399999999
$ws.Range("N64").Value = -4643.077
$ws.Range("H67").Value = 322558.78
$ws.Range("I67").Value = 540419.4399999999
$ws.Range("J67").Value = 4147.077
$ws.Range("K67").Value = 540419.4399999999
$ws.Range("L67").Value = 4147.077
$ws.Range("M67").Value = -539561.4399999999
$ws.Range("N67").Value = -5863.077
$ws.Range("H80").Value = 5444.6816
$ws.Range("I80").Value = 408.14285
$ws.Range("J80").Value = 14258.625
$ws.Range("K80").Value = 1224.42855
$ws.Range("L80").Value = 42775.875
$ws.Range("M80").Value = -226.4285500000001
$ws.Range("N80").Value = -44771.875
$ws.Range("H83").Value = 5444.6816
$ws.Range("I83").Value = 408.14285
$ws.Range("J83").Value = 14258.625
$ws.Range("K83").Value = 3673.28565
$ws.Range("L83").Value = 128327.625
$ws.Range("M83").Value = 1318.71435
$ws.Range("N83").Value = -138311.625
$ws.Range("H100").Value = 1291.8
$ws.Range("I100").Value = 1141.091
$ws.Range("J100").Value = 1706.25
$ws.Range("K100").Value = 1141.091
$ws.Range("L100").Value = 1706.25
$ws.Range("M100").Value = -600.0909999999999
$ws.Range("N100").Value = -2788.25
$ws.Range("H107").Value = 1065.3334
$ws.Range("I107").Value = 868.2857
$ws.Range("J107").Value = 1341.2
$ws.Range("K107").Value = 868.2857
$ws.Range("L107").Value = 1341.2
$ws.Range("M107").Value = 1051.7143
$ws.Range("N107").Value = -5181.2
$ws.Range("H137").Value = 3176.0715
$ws.Range("I137").Value = 2480.7144
$ws.Range("J137").Value = 3871.4285
$ws.Range("K137").Value = 7442.1432
$ws.Range("L137").Value = 11614.2855
$ws.Range("M137").Value = -4892.1432
$ws.Range("N137").Value = -16714.2855
$ws.Range("H138").Value = 4276240.5
$ws.Range("I138").Value = 1508
$ws.Range("J138").Value = 7411044.5
$ws.Range("K138").Value = 4524
$ws.Range("L138").Value = 22233133.5
$ws.Range("M138").Value = 616
$ws.Range("N138").Value = -22243413.5
$ws.Range("H141").Value = 5512.7554
$ws.Range("I141").Value = 2329.8276
$ws.Range("J141").Value = 11281.8125
$ws.Range("K141").Value = 6989.4828
$ws.Range("L141").Value = 33845.4375
$ws.Range("M141").Value = -1809.4828
$ws.Range("N141").Value = -44205.4375

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1872.2858
$ws.Range("I45").Value = 1820
$ws.Range("J45").Value = 2003
$ws.Range("K45").Value = 1820
$ws.Range("L45").Value = 2003
$ws.Range("M45").Value = -1443
$ws.Range("N45").Value = -2757
$ws.Range("H122").Value = 4311511.5
$ws.Range("I122").Value = 1155
$ws.Range("J122").Value = 15626198
$ws.Range("K122").Value = 3465
$ws.Range("L122").Value = 46878594
$ws.Range("M122").Value = -1015
$ws.Range("N122").Value = -46883494

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1498.3636
$ws.Range("I107").Value = 997.4286
$ws.Range("J107").Value = 2375
$ws.Range("K107").Value = 997.4286
$ws.Range("L107").Value = 2375
$ws.Range("M107").Value = 922.5714
$ws.Range("N107").Value = -6215
$ws.Range("H123").Value = 77032.5
$ws.Range("J123").Value = 77032.5
$ws.Range("L123").Value = 77032.5
$ws.Range("N123").Value = -86832.5
$ws.Range("H134").Value = 41791.08
$ws.Range("I134").Value = 3225.2222
$ws.Range("J134").Value = 128564.25
$ws.Range("K134").Value = 9675.6666
$ws.Range("L134").Value = 385692.75
$ws.Range("M134").Value = -7140.6666
$ws.Range("N134").Value = -390762.75

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1957.0652
$ws.Range("I31").Value = 1199.9706
$ws.Range("J31").Value = 4102.1665
$ws.Range("K31").Value = 1199.9706
$ws.Range("L31").Value = 4102.1665
$ws.Range("M31").Value = -904.9706000000001
$ws.Range("N31").Value = -4692.1665
$ws.Range("H34").Value = 1957.0652
$ws.Range("I34").Value = 1199.9706
$ws.Range("J34").Value = 4102.1665
$ws.Range("K34").Value = 1199.9706
$ws.Range("L34").Value = 4102.1665
$ws.Range("M34").Value = -997.9706000000001
$ws.Range("N34").Value = -4506.1665
$ws.Range("H39").Value = 4971.8335
$ws.Range("I39").Value = 4971.8335
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 4971.8335
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -4580.8335
$ws.Range("N39").ClearContents()
$ws.Range("H49").Value = 4971.8335
$ws.Range("I49").Value = 4971.8335
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 4971.8335
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -4789.8335
$ws.Range("N49").ClearContents()
$ws.Range("H50").Value = 15449.167
$ws.Range("J50").Value = 15449.167
$ws.Range("L50").Value = 15449.167
$ws.Range("N50").Value = -16699.167
$ws.Range("H58").Value = 1569586.9
$ws.Range("I58").Value = 2526825
$ws.Range("J58").Value = 3197.2727
$ws.Range("K58").Value = 2526825
$ws.Range("L58").Value = 3197.2727
$ws.Range("M58").Value = -2526622
$ws.Range("N58").Value = -3603.2727
$ws.Range("H59").Value = 29000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 29000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 29000
$ws.Range("N59").Value = -31290
$ws.Range("M59").ClearContents()
$ws.Range("H68").Value = 40295
$ws.Range("J68").Value = 40295
$ws.Range("L68").Value = 40295
$ws.Range("N68").Value = -41793
$ws.Range("H71").Value = 40295
$ws.Range("J71").Value = 40295
$ws.Range("L71").Value = 120885
$ws.Range("N71").Value = -128373
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74", "N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77", "N77").ClearContents()
$ws.Range("H99").Value = 2104.5715
$ws.Range("I99").Value = 1505.4667
$ws.Range("J99").Value = 3602.3333
$ws.Range("K99").Value = 1505.4667
$ws.Range("L99").Value = 3602.3333
$ws.Range("M99").Value = -7.466699999999946
$ws.Range("N99").Value = -6598.3333
$ws.Range("H122").Value = 7633.3706
$ws.Range("I122").Value = 3290.7646
$ws.Range("K122").Value = 9872.293799999999
$ws.Range("M122").Value = -7422.293799999999
$ws.Range("H126").Value = 2104.5715
$ws.Range("I126").Value = 1505.4667
$ws.Range("J126").Value = 3602.3333
$ws.Range("K126").Value = 4516.4001
$ws.Range("L126").Value = 10806.9999
$ws.Range("M126").Value = -2046.4001
$ws.Range("N126").Value = -15746.9999
$ws.Range("H134").Value = 2591.1936
$ws.Range("I134").Value = 1749.2903
$ws.Range("J134").Value = 3433.0967
$ws.Range("K134").Value = 5247.8709
$ws.Range("L134").Value = 10299.2901
$ws.Range("M134").Value = -2712.8709
$ws.Range("N134").Value = -15369.2901
$ws.Range("H136").Value = 1569586.9
$ws.Range("I136").Value = 2526825
$ws.Range("J136").Value = 3197.2727
$ws.Range("K136").Value = 7580475
$ws.Range("L136").Value = 9591.8181
$ws.Range("M136").Value = -7577925
$ws.Range("N136").Value = -14691.8181

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 567.9167
$ws.Range("I86").Value = 594.7143
$ws.Range("J86").Value = 530.4
$ws.Range("K86").Value = 1784.1429
$ws.Range("L86").Value = 1591.2
$ws.Range("M86").Value = -598.1428999999998
$ws.Range("N86").Value = -3963.2
$ws.Range("H89").Value = 567.9167
$ws.Range("I89").Value = 594.7143
$ws.Range("J89").Value = 530.4
$ws.Range("K89").Value = 5352.428699999999
$ws.Range("L89").Value = 4773.599999999999
$ws.Range("M89").Value = 575.5713000000005
$ws.Range("N89").Value = -16629.6
$ws.Range("H92").Value = 873.875
$ws.Range("J92").Value = 919.8
$ws.Range("L92").Value = 2759.4
$ws.Range("N92").Value = -5255.4
$ws.Range("H120").Value = 7732.1177
$ws.Range("I120").Value = 5405.8
$ws.Range("J120").Value = 8701.416999999999
$ws.Range("K120").Value = 16217.4
$ws.Range("L120").Value = 26104.251
$ws.Range("M120").Value = -11379.4
$ws.Range("N120").Value = -35780.251
$ws.Range("H121").Value = 1058
$ws.Range("J121").Value = 1162
$ws.Range("L121").Value = 3486
$ws.Range("N121").Value = -6106

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1193.3077
$ws.Range("I97").Value = 821.03705
$ws.Range("J97").Value = 2030.9166
$ws.Range("K97").Value = 821.03705
$ws.Range("L97").Value = 2030.9166
$ws.Range("M97").Value = -325.03705
$ws.Range("N97").Value = -3022.9166

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 563920.75
$ws.Range("I61").Value = 11146.182
$ws.Range("J61").Value = 1432566.4
$ws.Range("K61").Value = 11146.182
$ws.Range("L61").Value = 1432566.4
$ws.Range("M61").Value = -10944.182
$ws.Range("N61").Value = -1432970.4
$ws.Range("H64").Value = 33766.668
$ws.Range("J64").Value = 33766.668
$ws.Range("L64").Value = 33766.668
$ws.Range("N64").Value = -34216.668
$ws.Range("H67").Value = 33766.668
$ws.Range("J67").Value = 33766.668
$ws.Range("L67").Value = 33766.668
$ws.Range("N67").Value = -35326.668
$ws.Range("H93").Value = 424.10345
$ws.Range("I93").Value = 395.05264
$ws.Range("K93").Value = 395.05264
$ws.Range("M93").Value = 852.94736
$ws.Range("H113").Value = 563920.75
$ws.Range("I113").Value = 11146.182
$ws.Range("J113").Value = 1432566.4
$ws.Range("K113").Value = 11146.182
$ws.Range("L113").Value = 1432566.4
$ws.Range("M113").Value = -8976.182000000001
$ws.Range("N113").Value = -1436906.4
$ws.Range("H132").Value = 4371.3257
$ws.Range("I132").Value = 4085.4
$ws.Range("J132").Value = 5622.25
$ws.Range("K132").Value = 12256.2
$ws.Range("L132").Value = 16866.75
$ws.Range("M132").Value = -9726.200000000001
$ws.Range("N132").Value = -21926.75
$ws.Range("H136").Value = 4597.46
$ws.Range("I136").Value = 2800.25
$ws.Range("J136").Value = 7792.5
$ws.Range("K136").Value = 8400.75
$ws.Range("L136").Value = 23377.5
$ws.Range("M136").Value = -5850.75
$ws.Range("N136").Value = -28477.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3210.5
$ws.Range("I62").Value = 2918.0908
$ws.Range("K62").Value = 2918.0908
$ws.Range("M62").Value = -2294.0908
$ws.Range("H63").Value = 40249
$ws.Range("J63").Value = 40249
$ws.Range("L63").Value = 40249
$ws.Range("N63").Value = -41497
$ws.Range("H65").Value = 3210.5
$ws.Range("I65").Value = 2918.0908
$ws.Range("K65").Value = 14590.454
$ws.Range("M65").Value = -11470.454
$ws.Range("H66").Value = 40249
$ws.Range("J66").Value = 40249
$ws.Range("L66").Value = 120747
$ws.Range("N66").Value = -126987
$ws.Range("H96").Value = 1481.25
$ws.Range("J96").Value = 1410.6
$ws.Range("L96").Value = 1410.6
$ws.Range("N96").Value = -4156.6
$ws.Range("H132").Value = 1299.8474
$ws.Range("I132").Value = 502.2927
$ws.Range("J132").Value = 3116.5
$ws.Range("K132").Value = 1506.8781
$ws.Range("L132").Value = 9349.5
$ws.Range("M132").Value = 1023.1219
$ws.Range("N132").Value = -14409.5
$ws.Range("H136").Value = 4227.047
$ws.Range("I136").Value = 3249.8845
$ws.Range("J136").Value = 5766.8184
$ws.Range("K136").Value = 9749.6535
$ws.Range("L136").Value = 17300.4552
$ws.Range("M136").Value = -7199.6535
$ws.Range("N136").Value = -22400.4552
